$wb = $excel.ActiveWorkbook

# --- Rename column-B header labels on the two existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet as the last tab in the workbook ---
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Match the page margins used throughout the rest of the workbook
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Reuse the existing bold/centered/bordered header style for the new header row
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Reuse the existing date-number-format style for the new date column (A2:A68)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A68").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header row values ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---

$wsForecast.Range("A2").Value = 44976.99999999999
$wsForecast.Range("B2").Value = 154
$wsForecast.Range("C2").Value = -307.6317786613376
$wsForecast.Range("D2").Value = 614.6208507179261
$wsForecast.Range("A3").Value = 44983.99999999999
$wsForecast.Range("B3").Value = 159
$wsForecast.Range("C3").Value = -324.8199284071607
$wsForecast.Range("D3").Value = 616.3245792212094
$wsForecast.Range("A4").Value = 44997.99999999999
$wsForecast.Range("B4").Value = 168
$wsForecast.Range("C4").Value = -285.7711598263808
$wsForecast.Range("D4").Value = 657.0541769353091
$wsForecast.Range("A5").Value = 45011.99999999999
$wsForecast.Range("B5").Value = 178
$wsForecast.Range("C5").Value = -281.5991480315519
$wsForecast.Range("D5").Value = 662.7780936634776
$wsForecast.Range("A6").Value = 45025.99999999999
$wsForecast.Range("B6").Value = 188
$wsForecast.Range("C6").Value = -291.1984012092562
$wsForecast.Range("D6").Value = 656.0640520271767
$wsForecast.Range("A7").Value = 45032.99999999999
$wsForecast.Range("B7").Value = 192
$wsForecast.Range("C7").Value = -280.0065531933876
$wsForecast.Range("D7").Value = 688.3929106599411
$wsForecast.Range("A8").Value = 45053.99999999999
$wsForecast.Range("B8").Value = 207
$wsForecast.Range("C8").Value = -280.8912635298352
$wsForecast.Range("D8").Value = 687.5326018692929
$wsForecast.Range("A9").Value = 45060.99999999999
$wsForecast.Range("B9").Value = 212
$wsForecast.Range("C9").Value = -256.0140698252814
$wsForecast.Range("D9").Value = 699.388914229932
$wsForecast.Range("A10").Value = 45088.99999999999
$wsForecast.Range("B10").Value = 231
$wsForecast.Range("C10").Value = -202.7507455468546
$wsForecast.Range("D10").Value = 745.3018925373937
$wsForecast.Range("A11").Value = 45095.99999999999
$wsForecast.Range("B11").Value = 236
$wsForecast.Range("C11").Value = -249.6510573332441
$wsForecast.Range("D11").Value = 735.1291871007775
$wsForecast.Range("A12").Value = 45102.99999999999
$wsForecast.Range("B12").Value = 240
$wsForecast.Range("C12").Value = -224.7949400050063
$wsForecast.Range("D12").Value = 714.636482657309
$wsForecast.Range("A13").Value = 45116.99999999999
$wsForecast.Range("B13").Value = 250
$wsForecast.Range("C13").Value = -244.6957528136908
$wsForecast.Range("D13").Value = 747.0636752952387
$wsForecast.Range("A14").Value = 45123.99999999999
$wsForecast.Range("B14").Value = 255
$wsForecast.Range("C14").Value = -205.9105929661849
$wsForecast.Range("D14").Value = 745.148101944216
$wsForecast.Range("A15").Value = 45130.99999999999
$wsForecast.Range("B15").Value = 260
$wsForecast.Range("C15").Value = -184.1269909383631
$wsForecast.Range("D15").Value = 750.0166850473173
$wsForecast.Range("A16").Value = 45144.99999999999
$wsForecast.Range("B16").Value = 269
$wsForecast.Range("C16").Value = -188.4939549822456
$wsForecast.Range("D16").Value = 736.129479197756
$wsForecast.Range("A17").Value = 45151.99999999999
$wsForecast.Range("B17").Value = 274
$wsForecast.Range("C17").Value = -220.2297587215434
$wsForecast.Range("D17").Value = 728.2063969976897
$wsForecast.Range("A18").Value = 45158.99999999999
$wsForecast.Range("B18").Value = 279
$wsForecast.Range("C18").Value = -181.5582614264422
$wsForecast.Range("D18").Value = 759.463179729455
$wsForecast.Range("A19").Value = 45165.99999999999
$wsForecast.Range("B19").Value = 284
$wsForecast.Range("C19").Value = -217.8398803181302
$wsForecast.Range("D19").Value = 729.7164886794399
$wsForecast.Range("A20").Value = 45172.99999999999
$wsForecast.Range("B20").Value = 288
$wsForecast.Range("C20").Value = -155.241017450374
$wsForecast.Range("D20").Value = 759.4384897379807
$wsForecast.Range("A21").Value = 45179.99999999999
$wsForecast.Range("B21").Value = 293
$wsForecast.Range("C21").Value = -175.4786959651189
$wsForecast.Range("D21").Value = 797.7728462501691
$wsForecast.Range("A22").Value = 45186.99999999999
$wsForecast.Range("B22").Value = 298
$wsForecast.Range("C22").Value = -175.0233672855549
$wsForecast.Range("D22").Value = 798.7130400615674
$wsForecast.Range("A23").Value = 45193.99999999999
$wsForecast.Range("B23").Value = 303
$wsForecast.Range("C23").Value = -189.9070278182543
$wsForecast.Range("D23").Value = 756.6479534267662
$wsForecast.Range("A24").Value = 45200.99999999999
$wsForecast.Range("B24").Value = 308
$wsForecast.Range("C24").Value = -156.9094065800032
$wsForecast.Range("D24").Value = 809.0170793125026
$wsForecast.Range("A25").Value = 45221.99999999999
$wsForecast.Range("B25").Value = 322
$wsForecast.Range("C25").Value = -146.0777521114738
$wsForecast.Range("D25").Value = 823.27148365523
$wsForecast.Range("A26").Value = 45228.99999999999
$wsForecast.Range("B26").Value = 327
$wsForecast.Range("C26").Value = -143.7344523778972
$wsForecast.Range("D26").Value = 778.5323540471561
$wsForecast.Range("A27").Value = 45235.99999999999
$wsForecast.Range("B27").Value = 332
$wsForecast.Range("C27").Value = -159.4998142784946
$wsForecast.Range("D27").Value = 785.5168317531893
$wsForecast.Range("A28").Value = 45249.99999999999
$wsForecast.Range("B28").Value = 341
$wsForecast.Range("C28").Value = -107.0807128157744
$wsForecast.Range("D28").Value = 848.0405779237146
$wsForecast.Range("A29").Value = 45256.99999999999
$wsForecast.Range("B29").Value = 346
$wsForecast.Range("C29").Value = -132.5888038907824
$wsForecast.Range("D29").Value = 811.2461382818977
$wsForecast.Range("A30").Value = 45263.99999999999
$wsForecast.Range("B30").Value = 351
$wsForecast.Range("C30").Value = -145.3131733971231
$wsForecast.Range("D30").Value = 863.8861507066955
$wsForecast.Range("A31").Value = 45277.99999999999
$wsForecast.Range("B31").Value = 361
$wsForecast.Range("C31").Value = -131.888055328095
$wsForecast.Range("D31").Value = 831.4800921571848
$wsForecast.Range("A32").Value = 45298.99999999999
$wsForecast.Range("B32").Value = 375
$wsForecast.Range("C32").Value = -105.8105343273937
$wsForecast.Range("D32").Value = 838.4770521231502
$wsForecast.Range("A33").Value = 45305.99999999999
$wsForecast.Range("B33").Value = 380
$wsForecast.Range("C33").Value = -94.71577569961227
$wsForecast.Range("D33").Value = 866.2982040470112
$wsForecast.Range("A34").Value = 45312.99999999999
$wsForecast.Range("B34").Value = 385
$wsForecast.Range("C34").Value = -113.3813332021801
$wsForecast.Range("D34").Value = 839.7299746699076
$wsForecast.Range("A35").Value = 45326.99999999999
$wsForecast.Range("B35").Value = 394
$wsForecast.Range("C35").Value = -161.4692444042526
$wsForecast.Range("D35").Value = 842.6690682895467
$wsForecast.Range("A36").Value = 45333.99999999999
$wsForecast.Range("B36").Value = 399
$wsForecast.Range("C36").Value = -73.22635859821651
$wsForecast.Range("D36").Value = 904.5185966125491
$wsForecast.Range("A37").Value = 45340.99999999999
$wsForecast.Range("B37").Value = 404
$wsForecast.Range("C37").Value = -33.61265107788842
$wsForecast.Range("D37").Value = 895.6132436601106
$wsForecast.Range("A38").Value = 45347.99999999999
$wsForecast.Range("B38").Value = 409
$wsForecast.Range("C38").Value = -48.74842134315615
$wsForecast.Range("D38").Value = 890.9841990102509
$wsForecast.Range("A39").Value = 45354.99999999999
$wsForecast.Range("B39").Value = 413
$wsForecast.Range("C39").Value = -58.84231293439289
$wsForecast.Range("D39").Value = 927.1807658177692
$wsForecast.Range("A40").Value = 45361.99999999999
$wsForecast.Range("B40").Value = 418
$wsForecast.Range("C40").Value = -56.41827941456726
$wsForecast.Range("D40").Value = 888.7652519153576
$wsForecast.Range("A41").Value = 45375.99999999999
$wsForecast.Range("B41").Value = 428
$wsForecast.Range("C41").Value = -48.71167876665694
$wsForecast.Range("D41").Value = 872.6181559095315
$wsForecast.Range("A42").Value = 45382.99999999999
$wsForecast.Range("B42").Value = 433
$wsForecast.Range("C42").Value = -32.85853297964135
$wsForecast.Range("D42").Value = 886.3636629289316
$wsForecast.Range("A43").Value = 45389.99999999999
$wsForecast.Range("B43").Value = 437
$wsForecast.Range("C43").Value = -59.16707477202213
$wsForecast.Range("D43").Value = 896.0348118708165
$wsForecast.Range("A44").Value = 45396.99999999999
$wsForecast.Range("B44").Value = 442
$wsForecast.Range("C44").Value = -31.90336266348543
$wsForecast.Range("D44").Value = 933.7480131131069
$wsForecast.Range("A45").Value = 45403.99999999999
$wsForecast.Range("B45").Value = 447
$wsForecast.Range("C45").Value = -10.11329794308918
$wsForecast.Range("D45").Value = 916.6430421432136
$wsForecast.Range("A46").Value = 45417.99999999999
$wsForecast.Range("B46").Value = 457
$wsForecast.Range("C46").Value = 0.3400371109351905
$wsForecast.Range("D46").Value = 962.7504065253196
$wsForecast.Range("A47").Value = 45424.99999999999
$wsForecast.Range("B47").Value = 462
$wsForecast.Range("C47").Value = -33.17573455067121
$wsForecast.Range("D47").Value = 948.1470489573103
$wsForecast.Range("A48").Value = 45452.99999999999
$wsForecast.Range("B48").Value = 481
$wsForecast.Range("C48").Value = -6.94396505153827
$wsForecast.Range("D48").Value = 974.3712800913622
$wsForecast.Range("A49").Value = 45459.99999999999
$wsForecast.Range("B49").Value = 486
$wsForecast.Range("C49").Value = 48.0224806764952
$wsForecast.Range("D49").Value = 916.5875637312597
$wsForecast.Range("A50").Value = 45466.99999999999
$wsForecast.Range("B50").Value = 490
$wsForecast.Range("C50").Value = -12.30069027153389
$wsForecast.Range("D50").Value = 1029.587207045252
$wsForecast.Range("A51").Value = 45515.99999999999
$wsForecast.Range("B51").Value = 524
$wsForecast.Range("C51").Value = 48.35625987863333
$wsForecast.Range("D51").Value = 1012.424653703171
$wsForecast.Range("A52").Value = 45522.99999999999
$wsForecast.Range("B52").Value = 529
$wsForecast.Range("C52").Value = 70.03219134495856
$wsForecast.Range("D52").Value = 1023.849498045738
$wsForecast.Range("A53").Value = 45529.99999999999
$wsForecast.Range("B53").Value = 534
$wsForecast.Range("C53").Value = 35.3820554003057
$wsForecast.Range("D53").Value = 994.7155704980845
$wsForecast.Range("A54").Value = 45536.99999999999
$wsForecast.Range("B54").Value = 538
$wsForecast.Range("C54").Value = 74.112820629418
$wsForecast.Range("D54").Value = 1025.807853044468
$wsForecast.Range("A55").Value = 45557.99999999999
$wsForecast.Range("B55").Value = 553
$wsForecast.Range("C55").Value = 67.71155044416538
$wsForecast.Range("D55").Value = 1006.381008845506
$wsForecast.Range("A56").Value = 45564.99999999999
$wsForecast.Range("B56").Value = 558
$wsForecast.Range("C56").Value = 79.50579720820311
$wsForecast.Range("D56").Value = 1060.415855795049
$wsForecast.Range("A57").Value = 45592.99999999999
$wsForecast.Range("B57").Value = 577
$wsForecast.Range("C57").Value = 82.5560999044541
$wsForecast.Range("D57").Value = 1068.333742726902
$wsForecast.Range("A58").Value = 45599.99999999999
$wsForecast.Range("B58").Value = 582
$wsForecast.Range("C58").Value = 117.4345859825931
$wsForecast.Range("D58").Value = 1082.295814596984
$wsForecast.Range("A59").Value = 45634.99999999999
$wsForecast.Range("B59").Value = 606
$wsForecast.Range("C59").Value = 125.7206145063596
$wsForecast.Range("D59").Value = 1066.159986715952
$wsForecast.Range("A60").Value = 45641.99999999999
$wsForecast.Range("B60").Value = 611
$wsForecast.Range("C60").Value = 141.9730907755709
$wsForecast.Range("D60").Value = 1098.524899007873
$wsForecast.Range("A61").Value = 45648.99999999999
$wsForecast.Range("B61").Value = 615
$wsForecast.Range("C61").Value = 161.870358678308
$wsForecast.Range("D61").Value = 1104.568172498122
$wsForecast.Range("A62").Value = 45655.99999999999
$wsForecast.Range("B62").Value = 620
$wsForecast.Range("C62").Value = 120.6005581545765
$wsForecast.Range("D62").Value = 1095.099528858809
$wsForecast.Range("A63").Value = 45662.99999999999
$wsForecast.Range("B63").Value = 625
$wsForecast.Range("C63").Value = 132.4060711036983
$wsForecast.Range("D63").Value = 1082.844002994072
$wsForecast.Range("A64").Value = 45669.99999999999
$wsForecast.Range("B64").Value = 630
$wsForecast.Range("C64").Value = 122.1811282964439
$wsForecast.Range("D64").Value = 1137.882087152921
$wsForecast.Range("A65").Value = 45676.99999999999
$wsForecast.Range("B65").Value = 635
$wsForecast.Range("C65").Value = 147.2141101886543
$wsForecast.Range("D65").Value = 1079.559038431925
$wsForecast.Range("A66").Value = 45683.99999999999
$wsForecast.Range("B66").Value = 639
$wsForecast.Range("C66").Value = 144.1705640783291
$wsForecast.Range("D66").Value = 1104.901353118275
$wsForecast.Range("A67").Value = 45690.99999999999
$wsForecast.Range("B67").Value = 644
$wsForecast.Range("C67").Value = 180.2800366100052
$wsForecast.Range("D67").Value = 1121.51467586858
$wsForecast.Range("A68").Value = 45697.99999999999
$wsForecast.Range("B68").Value = 649
$wsForecast.Range("C68").Value = 159.9501711770882
$wsForecast.Range("D68").Value = 1086.354280733103

# --- Restore original active sheet/selection (first tab), matching the source workbook ---
$wsWeekly.Activate()
